# Optuna Attempt (go back with original)
# Update forecast values on "Forecast Comparison" and recalculated
# summary metrics on "Summary".

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------

# Row 2 (W8)
$wsForecast.Range("D2").Value = 6
$wsForecast.Range("H2").Value = 20.17
$wsForecast.Range("L2").Value = 0.8100000000000001

# Row 3 (W9)
$wsForecast.Range("D3").Value = 7
$wsForecast.Range("H3").Value = 18.08
$wsForecast.Range("L3").Value = 0.86

# Row 4 (W10)
$wsForecast.Range("D4").Value = 7
$wsForecast.Range("H4").Value = 17.22
$wsForecast.Range("L4").Value = 0.86

# Row 5 (W11)
$wsForecast.Range("H5").Value = 13.13
$wsForecast.Range("L5").Value = 1.03

# Row 6 (W12)
$wsForecast.Range("H6").Value = 11.79
$wsForecast.Range("L6").Value = 1.01

# Row 7 (W13)
$wsForecast.Range("H7").Value = 11.1
$wsForecast.Range("L7").Value = 1.11

# Row 8 (W14)
$wsForecast.Range("H8").Value = 9.81
$wsForecast.Range("L8").Value = 1.07

# Row 9 (W15)
$wsForecast.Range("D9").Value = 6
$wsForecast.Range("H9").Value = 12.21
$wsForecast.Range("L9").Value = 0.96

# Row 10 (W16)
$wsForecast.Range("D10").Value = 7
$wsForecast.Range("H10").Value = 10.01
$wsForecast.Range("L10").Value = 1.04

# Row 11 (W17)
$wsForecast.Range("D11").Value = 8
$wsForecast.Range("H11").Value = 8.140000000000001
$wsForecast.Range("L11").Value = 1

# Row 12 (W18)
$wsForecast.Range("H12").Value = 6.38
$wsForecast.Range("L12").Value = 0.82

# Row 13 (W19)
$wsForecast.Range("H13").Value = 5.24
$wsForecast.Range("L13").Value = 0.82

# Row 14 (W20)
$wsForecast.Range("H14").Value = 4.36
$wsForecast.Range("L14").Value = 0.93

# Row 15 (W21)
$wsForecast.Range("H15").Value = 3.46
$wsForecast.Range("L15").Value = 0.8100000000000001

# Row 16 (W22)
$wsForecast.Range("H16").Value = 2.46
$wsForecast.Range("L16").Value = 0.85

# Row 17 (W23)
$wsForecast.Range("H17").Value = 1.46
$wsForecast.Range("I17").Value = "Low"
$wsForecast.Range("J17").Value = "Normal"
$wsForecast.Range("L17").Value = 1.02

# --- Summary sheet ---------------------------------------------------------
# These "Value" cells hold numeric-looking text (t="inlineStr" in the
# source), so force text formatting before assigning, matching how Excel
# keeps a numeric-looking entry as text.

$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "136"

$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value = "65"

$wsSummary.Range("B11").NumberFormat = "@"
$wsSummary.Range("B11").Value = "31"

$wsSummary.Range("B14").NumberFormat = "@"
$wsSummary.Range("B14").Value = "7"
